$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @{ "B"=14.97980646600029; "C"=11.54169036529256; "D"=5.956003378244518; "E"=16.60567352243067; "G"=18.88952722893729; "H"=11.55325594156004; "I"=15.41205661939161; "O"=16.27990738049906 }
    3 = @{ "B"=14.12598843277608; "C"=10.88308628679088; "D"=5.831707957588392; "E"=15.65631050192974; "G"=18.8376966261425; "H"=11.61431219392075; "I"=15.57718986152581; "O"=16.35947129485095 }
    4 = @{ "B"=13.57368922920541; "C"=10.45603610136106; "D"=5.755835246328446; "E"=15.04796707570479; "G"=18.81980491406887; "H"=11.65499083769751; "I"=15.68396404223561; "O"=16.41512930587027 }
    5 = @{ "B"=13.34170703070008; "C"=10.27638548787544; "D"=5.725077168354272; "E"=14.79393393156256; "G"=18.81600375420437; "H"=11.67236673586066; "I"=15.72882890540635; "O"=16.43950734213673 }
    6 = @{ "B"=13.30277351906185; "C"=10.24621750930319; "D"=5.719980953679795; "E"=14.75139062807681; "G"=18.81558279203027; "H"=11.67530015927673; "I"=15.73636044107148; "O"=16.44365736114525 }
    7 = @{ "B"=13.5705884327998; "C"=10.45363593000209; "D"=5.755419716705041; "E"=15.04456551773486; "G"=18.81973954427913; "H"=11.65522194405117; "I"=15.68456362475538; "O"=16.41545122503712 }
    8 = @{ "B"=14.69130566510876; "C"=11.31934729866441; "D"=5.913082909871251; "E"=16.28375547595196; "G"=18.86875880224762; "H"=11.57364421738025; "I"=15.46787770527352; "O"=16.30592004385321 }
    9 = @{ "B"=16.66180294336104; "C"=12.83470874885098; "D"=6.223699370150761; "E"=18.62418747806831; "G"=19.0756070028165; "H"=11.43911964697766; "I"=15.08564995758072; "O"=16.14574923331975 }
    10 = @{ "B"=17.9663646206055; "C"=13.83454840564394; "D"=6.45008390290055; "E"=20.28156041951576; "G"=19.29469451258568; "H"=11.35599011114372; "I"=14.83084647940949; "O"=16.06218686317134 }
    11 = @{ "B"=18.52813304948098; "C"=14.26446411632242; "D"=6.552140911892568; "E"=20.99305485599736; "G"=19.40867321735213; "H"=11.32162414129091; "I"=14.72058228150835; "O"=16.03174367051663 }
    12 = @{ "B"=18.736268629856; "C"=14.42366383577119; "D"=6.590612555693451; "E"=21.25641184347996; "G"=19.45385823018928; "H"=11.30911025029682; "I"=14.67964106146328; "O"=16.021316641424 }
    13 = @{ "B"=18.69164771834321; "C"=14.38953764882966; "D"=6.582335428215632; "E"=21.19996265048718; "G"=19.44403745596546; "H"=11.31178305574764; "I"=14.68842228299624; "O"=16.02351312145804 }
    14 = @{ "B"=18.54534873826896; "C"=14.27763379961502; "D"=6.555309738123335; "E"=21.0148428611808; "G"=19.41235028250784; "H"=11.32058457920411; "I"=14.71719771132467; "O"=16.03086369415314 }
    15 = @{ "B"=18.45513717727153; "C"=14.20862022646098; "D"=6.538731714684499; "E"=20.90066209460528; "G"=19.39320336517947; "H"=11.32604095766738; "I"=14.7349294712571; "O"=16.03550989659608 }
    16 = @{ "B"=17.92901069744875; "C"=13.80594942030116; "D"=6.443391955198709; "E"=20.23421124982915; "G"=19.28753099696461; "H"=11.35830574320667; "I"=14.83816625184556; "O"=16.06432981557364 }
    17 = @{ "B"=17.59810383853685; "C"=13.5525280867111; "D"=6.384636206849326; "E"=19.81451642560857; "G"=19.22634856787156; "H"=11.37898537591142; "I"=14.90294595177865; "O"=16.08395790255538 }
    18 = @{ "B"=17.40479429438562; "C"=13.40442217552418; "D"=6.350756604899041; "E"=19.56912117526506; "G"=19.1925080092936; "H"=11.39120427213403; "I"=14.94073696264739; "O"=16.09595924692483 }
    19 = @{ "B"=17.33883240484055; "C"=13.35387398347001; "D"=6.339272298655141; "E"=19.48534687442998; "G"=19.18128296294844; "H"=11.39539700617249; "I"=14.95362358734572; "O"=16.1001445500298 }
    20 = @{ "B"=17.63363823519415; "C"=13.57974802717059; "D"=6.390899949730308; "E"=19.85960717454698; "G"=19.23272205302159; "H"=11.37675038090123; "I"=14.89599503076437; "O"=16.08179469803346 }
    21 = @{ "B"=18.58844522345726; "C"=14.31060048298277; "D"=6.563252912871799; "E"=21.06938151162281; "G"=19.42160296388493; "H"=11.31798576610932; "I"=14.70872358168376; "O"=16.02867467087193 }
    22 = @{ "B"=19.1856792324153; "C"=14.76726715814206; "D"=6.674857107322819; "E"=21.82467106298138; "G"=19.55682214890962; "H"=11.28249435020563; "I"=14.59107303930695; "O"=16.00038133104208 }
    23 = @{ "B"=18.86938517622873; "C"=14.52546007690663; "D"=6.615400052499676; "E"=21.42478375538209; "G"=19.48358916195222; "H"=11.30116889002256; "I"=14.65343091296374; "O"=16.01489026930843 }
    24 = @{ "B"=17.61758268385637; "C"=13.56744939937391; "D"=6.388068420297493; "E"=19.83923444302784; "G"=19.22983644232337; "H"=11.37775979473187; "I"=14.89913583421757; "O"=16.08277045056931 }
    25 = @{ "B"=16.153610737176; "C"=12.44456650955127; "D"=6.139801196593193; "E"=17.97597161059372; "G"=19.00780092786124; "H"=11.47276794407263; "I"=15.1844837104312; "O"=16.18314583140229 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value2 = $data[$row][$col]
    }
}
